# Atualizado por script em 27-11-2023 08:45
#
# 1) Rows 42/43 and 99/100: the match pairs (same kickoff date/time) were
#    recorded in the wrong order - swap the match-specific columns (F..V)
#    between each pair while leaving the index/metadata columns (A..E)
#    untouched.
# 2) Append three new match rows (221..223) at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-MatchColumns($rowA, $rowB) {
    for ($col = 6; $col -le 22; $col++) {
        $cellA = $ws.Cells.Item($rowA, $col)
        $cellB = $ws.Cells.Item($rowB, $col)
        $valueA = $cellA.Value2
        $valueB = $cellB.Value2
        $cellA.Value = $valueB
        $cellB.Value = $valueA
    }
}

Swap-MatchColumns 42 43
Swap-MatchColumns 99 100

function Set-MatchRow(
    $row,
    $indice, $pais, $torneio, $temporada, $dataPartida,
    $home, $homeGols, $away, $awayGols,
    $homeOpenOdds, $homeOpenDT, $homeCloseOdds, $homeCloseDT,
    $drawOpenOdds, $drawOpenDT, $drawCloseOdds, $drawCloseDT,
    $awayOpenOdds, $awayOpenDT, $awayCloseOdds, $awayCloseDT,
    $url
) {
    # Clone formatting (styles + number formats) from the previous row,
    # then overwrite the values - this reuses the existing style records
    # (bold/bordered index column, date-formatted column E) instead of
    # creating new duplicate styles.
    $srcRow = $row - 1
    $ws.Range("A" + $srcRow + ":V" + $srcRow).Copy($ws.Range("A" + $row + ":V" + $row))

    $ws.Cells.Item($row, 1).Value = $indice
    $ws.Cells.Item($row, 2).Value = $pais
    $ws.Cells.Item($row, 3).Value = $torneio

    # "temporada" is a numeric-looking string ("2023") that must stay text
    # (matches the other rows, which store it as inlineStr). Force the Text
    # number format before the write so it is not auto-parsed as a number,
    # then restore the default style so no stray formatting is left behind.
    $temporadaCell = $ws.Cells.Item($row, 4)
    $temporadaCell.NumberFormat = "@"
    $temporadaCell.Value = $temporada
    $temporadaCell.Style = "Normal"

    $ws.Cells.Item($row, 5).Value = $dataPartida
    $ws.Cells.Item($row, 6).Value = $home
    $ws.Cells.Item($row, 7).Value = $homeGols
    $ws.Cells.Item($row, 8).Value = $away
    $ws.Cells.Item($row, 9).Value = $awayGols
    $ws.Cells.Item($row, 10).Value = $homeOpenOdds
    $ws.Cells.Item($row, 11).Value = $homeOpenDT
    $ws.Cells.Item($row, 12).Value = $homeCloseOdds
    $ws.Cells.Item($row, 13).Value = $homeCloseDT
    $ws.Cells.Item($row, 14).Value = $drawOpenOdds
    $ws.Cells.Item($row, 15).Value = $drawOpenDT
    $ws.Cells.Item($row, 16).Value = $drawCloseOdds
    $ws.Cells.Item($row, 17).Value = $drawCloseDT
    $ws.Cells.Item($row, 18).Value = $awayOpenOdds
    $ws.Cells.Item($row, 19).Value = $awayOpenDT
    $ws.Cells.Item($row, 20).Value = $awayCloseOdds
    $ws.Cells.Item($row, 21).Value = $awayCloseDT
    $ws.Cells.Item($row, 22).Value = $url
}

Set-MatchRow 221 220 "chile" "primera-division" "2023" 45256.91666666666 `
    "A. Italiano" 0 "Colo Colo" 1 `
    3.33 "19/11/2023 03:42" 4.58 "26/11/2023 21:44" `
    3.44 "19/11/2023 03:42" 3.71 "26/11/2023 21:53" `
    2.21 "19/11/2023 03:42" 1.83 "26/11/2023 21:44" `
    "https://www.betexplorer.com/football/chile/primera-division/a-italiano-colo-colo/MmIPVoFb/"

Set-MatchRow 222 221 "chile" "primera-division" "2023" 45257.02083333334 `
    "Curico Unido" 3 "Magallanes" 4 `
    2.01 "19/11/2023 03:43" 2.66 "27/11/2023 00:30" `
    3.56 "19/11/2023 03:43" 3.63 "27/11/2023 00:00" `
    3.51 "19/11/2023 03:43" 2.62 "27/11/2023 00:30" `
    "https://www.betexplorer.com/football/chile/primera-division/curico-unido-magallanes/8pjDMsxo/"

Set-MatchRow 223 222 "chile" "primera-division" "2023" 45257.02083333334 `
    "Copiapo" 1 "Nublense" 1 `
    2.49 "19/11/2023 03:42" 2.8 "27/11/2023 00:26" `
    3.47 "19/11/2023 03:42" 3.2 "27/11/2023 00:26" `
    2.92 "19/11/2023 03:42" 2.76 "27/11/2023 00:26" `
    "https://www.betexplorer.com/football/chile/primera-division/copiapo-nublense/MmveoREo/"
